# ------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q3" sheet (it already has the correct column
#    layout/styling) and place the copy right after it, then rename
#    it to "2022-Q1". This keeps it positioned before "总计", matching
#    the target sheet order:
#       2020-Q4, 2021-Q2, 2021-Q3, 2022-Q1, 总计
# 2. Overwrite the header + data cells of the new sheet with the
#    2022-Q1 numbers.
# 3. Insert the new 2022-Q1 summary row at the top of "总计",
#    shifting the existing rows down by one.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q1" sheet from a copy of "2021-Q3" ---
$q3 = $wb.Worksheets.Item("2021-Q3")
$q3.Copy($null, $q3)
$newSheet = $wb.Worksheets.Item(4)
$newSheet.Name = "2022-Q1"

# --- Step 2: update header row (D1 label changed) ---
$newSheet.Range("D1").Value = "基金规模"

# Fund code (B) and amount/position/ratio/value (D,E,F,G) columns must
# stay text, otherwise Excel will coerce numeric-looking strings (e.g.
# the leading-zero fund codes, or trailing-zero decimals) into numbers.
# Force text format, assign the values, then restore the "Normal" style
# so no stray number-format style sticks around on the cells.
$newSheet.Range("B2:B8").NumberFormat = "@"
$newSheet.Range("D2:G8").NumberFormat = "@"

$newSheet.Range("B2").Value = "005613"
$newSheet.Range("C2").Value = "上投摩根富时发达市场REITs指数QDII人民币份额"
$newSheet.Range("D2").Value = "4.84"
$newSheet.Range("E2").Value = "91.10"
$newSheet.Range("F2").Value = "3.26"
$newSheet.Range("G2").Value = "0.1578"
$newSheet.Range("H2").Value = 6

$newSheet.Range("B3").Value = "005614"
$newSheet.Range("C3").Value = "上投摩根富时发达市场REITs指数QDII美钞"
$newSheet.Range("D3").Value = "4.84"
$newSheet.Range("E3").Value = "91.10"
$newSheet.Range("F3").Value = "3.26"
$newSheet.Range("G3").Value = "0.1578"
$newSheet.Range("H3").Value = 6

$newSheet.Range("B4").Value = "005615"
$newSheet.Range("C4").Value = "上投摩根富时发达市场REITs指数QDII美汇"
$newSheet.Range("D4").Value = "4.84"
$newSheet.Range("E4").Value = "91.10"
$newSheet.Range("F4").Value = "3.26"
$newSheet.Range("G4").Value = "0.1578"
$newSheet.Range("H4").Value = 6

$newSheet.Range("B5").Value = "000179"
$newSheet.Range("C5").Value = "广发美国房地产指数QDII-人民币"
$newSheet.Range("D5").Value = "2.37"
$newSheet.Range("E5").Value = "92.38"
$newSheet.Range("F5").Value = "2.62"
$newSheet.Range("G5").Value = "0.0621"
$newSheet.Range("H5").Value = 8

$newSheet.Range("B6").Value = "000180"
$newSheet.Range("C6").Value = "广发美国房地产指数QDII - 美元"
$newSheet.Range("D6").Value = "2.37"
$newSheet.Range("E6").Value = "92.38"
$newSheet.Range("F6").Value = "2.62"
$newSheet.Range("G6").Value = "0.0621"
$newSheet.Range("H6").Value = 8

$newSheet.Range("B7").Value = "160140"
$newSheet.Range("C7").Value = "南方道琼斯美国精选REIT指数(QDII-LOF)A"
$newSheet.Range("D7").Value = "1.35"
$newSheet.Range("E7").Value = "89.10"
$newSheet.Range("F7").Value = "2.78"
$newSheet.Range("G7").Value = "0.0375"
$newSheet.Range("H7").Value = 7

$newSheet.Range("B8").Value = "160141"
$newSheet.Range("C8").Value = "南方道琼斯美国精选REIT指数(QDII-LOF)C"
$newSheet.Range("D8").Value = "0.44"
$newSheet.Range("E8").Value = "89.10"
$newSheet.Range("F8").Value = "2.78"
$newSheet.Range("G8").Value = "0.0122"
$newSheet.Range("H8").Value = 7

# Clean the stray text-number-format style back off these cells so they
# match the plain (unstyled) data cells used elsewhere in the sheet.
$newSheet.Range("B2:B8").Style = "Normal"
$newSheet.Range("D2:G8").Style = "Normal"

# --- Step 3: update the "总计" (total) summary sheet ---
$totalSheet = $wb.Worksheets.Item("总计")

# Shift existing rows 2-4 down to 3-5 (copy preserves formatting/style),
# processing bottom-up so sources aren't clobbered before being read.
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

# Fix up the running index in column A for the shifted rows.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Write the new 2022-Q1 summary row at the top.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.65
